$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 360 (shifts existing rows 360-415 down to 361-416)
$ws.Rows.Item(360).Insert()

# Populate the new row 360 with the new weekly record
$ws.Cells.Item(360, 1).Value = 11
$ws.Cells.Item(360, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(360, 3).Value = "Bíobío"
$ws.Cells.Item(360, 4).Value = 45127
$ws.Cells.Item(360, 5).Value = 8
$ws.Cells.Item(360, 6).Value = 100112045
$ws.Cells.Item(360, 7).Value = "Zapallo"
$ws.Cells.Item(360, 8).Value = "Camote"
$ws.Cells.Item(360, 9).Value = "1a (guarda)"
$ws.Cells.Item(360, 10).Value = 1000
$ws.Cells.Item(360, 11).Value = 400
$ws.Cells.Item(360, 12).Value = 500
$ws.Cells.Item(360, 13).Value = 450
$ws.Cells.Item(360, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(360, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(360, 16).Value = 450
$ws.Cells.Item(360, 17).Value = 1
$ws.Cells.Item(360, 18).Value = "Hortaliza"
